$d = $word.ActiveDocument

$replacements = @(
    @("2025-02-28 Friday", "2025-03-01 Saturday"),
    @("361×5=", "150×7="),
    @("358×9=", "715×6="),
    @("382×6=", "518×7="),
    @("732×4=", "289×2="),
    @("577×8=", "456×9="),
    @("685×3=", "329×7="),
    @("219×7=", "326×9="),
    @("558×9=", "481×7="),
    @("195×8=", "398×9="),
    @("443×3=", "725×4="),
    @("737×6=", "168×8="),
    @("998×8=", "738×2="),
    @("164×2=", "827×9="),
    @("795×7=", "698×6="),
    @("536×4=", "287×6="),
    @("673×2=", "957×8="),
    @("246×9=", "265×9="),
    @("723×9=", "898×5="),
    @("112×7=", "606×8="),
    @("251×5=", "911×6="),
    @("149×2=", "466×6="),
    @("520×9=", "321×8="),
    @("746×7=", "498×4="),
    @("906×7=", "842×7="),
    @("118×2=", "856×8=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
